$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The e-mail address stored in B8 was missing the "@" ("user_email7ext.com").
# Fix the typo so it reads "user_email7@ext.com", matching the pattern used
# by every other row in the "email" column.
$ws.Range("B8").Value = "user_email7@ext.com"

# Every other address in column B carries a "mailto:" hyperlink pointing at
# itself; B8 never got one because the malformed address wasn't recognised.
# Add it now that the address is valid.
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:user_email7@ext.com")

# Adding the hyperlink re-applies Excel's built-in hyperlink formatting as a
# brand new style; put B8 back on the same "Lien hypertexte" cell style the
# rest of column B already uses so formatting stays consistent.
$ws.Range("B8").Style = "Lien hypertexte"

# Move the active selection to the cell that was just corrected.
[void]$ws.Range("B8").Select()
